# Auto-generated edit script applying the Phantom_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 674.5
$ws.Range("I12").Value = 661
$ws.Range("J12").Value = 850
$ws.Range("K12").Value = 661
$ws.Range("L12").Value = 850
$ws.Range("M12").Value = -491
$ws.Range("N12").Value = -1190
$ws.Range("H18").Value = 5280.2
$ws.Range("J18").Value = 1002
$ws.Range("L18").Value = 1002
$ws.Range("N18").Value = -1570
$ws.Range("H62").Value = 3733.6667
$ws.Range("I62").Value = 3600.5
$ws.Range("K62").Value = 3600.5
$ws.Range("M62").Value = -2976.5
$ws.Range("H65").Value = 3733.6667
$ws.Range("I65").Value = 3600.5
$ws.Range("K65").Value = 18002.5
$ws.Range("M65").Value = -14882.5
$ws.Range("H76").Value = 40002500
$ws.Range("I76").Value = 66667970
$ws.Range("K76").Value = 66667970
$ws.Range("M76").Value = -66667655
$ws.Range("H79").Value = 40002500
$ws.Range("I79").Value = 66667970
$ws.Range("K79").Value = 66667970
$ws.Range("M79").Value = -66666878
$ws.Range("H88").Value = 6665.6665
$ws.Range("I88").Value = 4999
$ws.Range("K88").Value = 4999
$ws.Range("M88").Value = -4593
$ws.Range("H91").Value = 6665.6665
$ws.Range("I91").Value = 4999
$ws.Range("K91").Value = 4999
$ws.Range("M91").Value = -3595
$ws.Range("H92").Value = 100359.1
$ws.Range("I92").Value = 100359.1
$ws.Range("K92").Value = 100359.1
$ws.Range("M92").Value = -99111.10000000001
$ws.Range("H106").Value = 21856.285
$ws.Range("I106").Value = 20165.916
$ws.Range("K106").Value = 20165.916
$ws.Range("M106").Value = -19534.916
$ws.Range("H107").Value = 1580.2858
$ws.Range("I107").Value = 1442.9
$ws.Range("J107").Value = 1923.75
$ws.Range("K107").Value = 1442.9
$ws.Range("L107").Value = 1923.75
$ws.Range("M107").Value = 477.0999999999999
$ws.Range("N107").Value = -5763.75
$ws.Range("H113").Value = 1005
$ws.Range("I113").Value = 1005
$ws.Range("K113").Value = 1005
$ws.Range("M113").Value = 2249
$ws.Range("H116").Value = 4925
$ws.Range("I116").Value = 4887.5
$ws.Range("K116").Value = 4887.5
$ws.Range("M116").Value = -1445.5
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = $null
$ws.Range("H132").Value = 5661.6665
$ws.Range("I132").Value = 5661.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16984.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14454.9995
$ws.Range("N132").Value = $null
$ws.Range("H137").Value = 2613.36
$ws.Range("I137").Value = 2268.9524
$ws.Range("J137").Value = 4421.5
$ws.Range("K137").Value = 6806.8572
$ws.Range("L137").Value = 13264.5
$ws.Range("M137").Value = -4256.8572
$ws.Range("N137").Value = -18364.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1303.3
$ws.Range("I2").Value = 772.3333
$ws.Range("J2").Value = 2099.75
$ws.Range("K2").Value = 772.3333
$ws.Range("L2").Value = 2099.75
$ws.Range("M2").Value = -659.3333
$ws.Range("N2").Value = -2325.75
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("H74").Value = 1880.5385
$ws.Range("I74").Value = 1829.0834
$ws.Range("K74").Value = 1829.0834
$ws.Range("M74").Value = -955.0834
$ws.Range("H77").Value = 1880.5385
$ws.Range("I77").Value = 1829.0834
$ws.Range("K77").Value = 9145.416999999999
$ws.Range("M77").Value = -4777.416999999999
$ws.Range("H97").Value = 983.7143
$ws.Range("J97").Value = 5500
$ws.Range("L97").Value = 5500
$ws.Range("N97").Value = -6492
$ws.Range("H102").Value = 2108.4443
$ws.Range("I102").Value = 2108.4443
$ws.Range("K102").Value = 2108.4443
$ws.Range("M102").Value = -486.4443000000001
$ws.Range("H110").Value = 6919.75
$ws.Range("I110").Value = 7765.4287
$ws.Range("K110").Value = 7765.4287
$ws.Range("M110").Value = -5720.4287
$ws.Range("H116").Value = 1303.3
$ws.Range("I116").Value = 772.3333
$ws.Range("J116").Value = 2099.75
$ws.Range("K116").Value = 772.3333
$ws.Range("L116").Value = 2099.75
$ws.Range("M116").Value = 1521.6667
$ws.Range("N116").Value = -6687.75
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1303.3
$ws.Range("I3").Value = 772.3333
$ws.Range("J3").Value = 2099.75
$ws.Range("K3").Value = 772.3333
$ws.Range("L3").Value = 2099.75
$ws.Range("M3").Value = -658.3333
$ws.Range("N3").Value = -2327.75
$ws.Range("H64").Value = 2988.2727
$ws.Range("I64").Value = 1244.25
$ws.Range("J64").Value = 3984.8572
$ws.Range("K64").Value = 1244.25
$ws.Range("L64").Value = 3984.8572
$ws.Range("M64").Value = -1019.25
$ws.Range("N64").Value = -4434.8572
$ws.Range("H67").Value = 2988.2727
$ws.Range("I67").Value = 1244.25
$ws.Range("J67").Value = 3984.8572
$ws.Range("K67").Value = 1244.25
$ws.Range("L67").Value = 3984.8572
$ws.Range("M67").Value = -464.25
$ws.Range("N67").Value = -5544.8572
$ws.Range("H86").Value = 14571.2
$ws.Range("I86").Value = 14571.2
$ws.Range("K86").Value = 14571.2
$ws.Range("M86").Value = -13448.2
$ws.Range("H89").Value = 14571.2
$ws.Range("I89").Value = 14571.2
$ws.Range("K89").Value = 72856
$ws.Range("M89").Value = -67240
$ws.Range("H94").Value = 359.14285
$ws.Range("I94").Value = 359.14285
$ws.Range("K94").Value = 359.14285
$ws.Range("M94").Value = 91.85714999999999
$ws.Range("H99").Value = 1999
$ws.Range("I99").Value = 1999
$ws.Range("K99").Value = 1999
$ws.Range("M99").Value = -501
$ws.Range("H105").Value = 4814.357
$ws.Range("I105").Value = 4296.4
$ws.Range("K105").Value = 4296.4
$ws.Range("M105").Value = -2549.4
$ws.Range("H107").Value = 1506
$ws.Range("I107").Value = 1506
$ws.Range("K107").Value = 1506
$ws.Range("M107").Value = 414
$ws.Range("H134").Value = 5964.533
$ws.Range("I134").Value = 6204.5
$ws.Range("K134").Value = 18613.5
$ws.Range("M134").Value = -16078.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2735.1428
$ws.Range("I16").Value = 2786.75
$ws.Range("K16").Value = 2786.75
$ws.Range("M16").Value = -2499.75
$ws.Range("H22").Value = 4444919.5
$ws.Range("I22").Value = 554.2
$ws.Range("J22").Value = 10000376
$ws.Range("K22").Value = 554.2
$ws.Range("L22").Value = 10000376
$ws.Range("M22").Value = -204.2
$ws.Range("N22").Value = -10001076
$ws.Range("H94").Value = 1475
$ws.Range("J94").Value = 1025
$ws.Range("L94").Value = 1025
$ws.Range("N94").Value = -1927
$ws.Range("H99").Value = 3071.2856
$ws.Range("I99").Value = 2999.8333
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2999.8333
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1501.8333
$ws.Range("N99").Value = -6496
$ws.Range("H105").Value = 2282.2
$ws.Range("I105").Value = 2026.375
$ws.Range("J105").Value = 3305.5
$ws.Range("K105").Value = 2026.375
$ws.Range("L105").Value = 3305.5
$ws.Range("M105").Value = -279.375
$ws.Range("N105").Value = -6799.5
$ws.Range("H113").Value = 2735.1428
$ws.Range("I113").Value = 2786.75
$ws.Range("K113").Value = 2786.75
$ws.Range("M113").Value = -616.75
$ws.Range("H126").Value = 3071.2856
$ws.Range("I126").Value = 2999.8333
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 8999.499899999999
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -6529.499899999999
$ws.Range("N126").Value = -15440
$ws.Range("H135").Value = 89999.5
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 100477.5
$ws.Range("I7").Value = 133783.67
$ws.Range("J7").Value = 559
$ws.Range("K7").Value = 401351.01
$ws.Range("L7").Value = 1677
$ws.Range("M7").Value = -401239.01
$ws.Range("N7").Value = -1901
$ws.Range("H11").Value = 747.75
$ws.Range("I11").Value = 727.7143
$ws.Range("K11").Value = 2183.1429
$ws.Range("M11").Value = -2043.1429
$ws.Range("H17").Value = 744
$ws.Range("I17").Value = 722.5
$ws.Range("K17").Value = 2167.5
$ws.Range("M17").Value = -1998.5
$ws.Range("H63").Value = 9997
$ws.Range("I63").Value = 9997
$ws.Range("K63").Value = 29991
$ws.Range("M63").Value = -29242
$ws.Range("H66").Value = 9997
$ws.Range("I66").Value = 9997
$ws.Range("K66").Value = 89973
$ws.Range("M66").Value = -86229
$ws.Range("H113").Value = 3076.3333
$ws.Range("J113").Value = 2992.4546
$ws.Range("L113").Value = 8977.363799999999
$ws.Range("N113").Value = -13317.3638
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 109.25
$ws.Range("I2").Value = 122
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 122
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = -9
$ws.Range("N2").Value = -246
$ws.Range("H70").Value = 24891.23
$ws.Range("I70").Value = 30158.7
$ws.Range("K70").Value = 30158.7
$ws.Range("M70").Value = -29888.7
$ws.Range("H73").Value = 24891.23
$ws.Range("I73").Value = 30158.7
$ws.Range("K73").Value = 30158.7
$ws.Range("M73").Value = -29222.7
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H108").Value = 66665.664
$ws.Range("J108").Value = 66665.664
$ws.Range("L108").Value = 66665.664
$ws.Range("N108").Value = -74345.664
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 460.63635
$ws.Range("I22").Value = 331.77777
$ws.Range("J22").Value = 1040.5
$ws.Range("K22").Value = 331.77777
$ws.Range("L22").Value = 1040.5
$ws.Range("M22").Value = -36.77776999999998
$ws.Range("N22").Value = -1630.5
$ws.Range("H27").Value = 460.63635
$ws.Range("I27").Value = 331.77777
$ws.Range("J27").Value = 1040.5
$ws.Range("K27").Value = 331.77777
$ws.Range("L27").Value = 1040.5
$ws.Range("M27").Value = -224.77777
$ws.Range("N27").Value = -1254.5
$ws.Range("H40").Value = 2318.6667
$ws.Range("I40").Value = 2318.6667
$ws.Range("K40").Value = 2318.6667
$ws.Range("M40").Value = -2182.6667
$ws.Range("H61").Value = 1748.1333
$ws.Range("I61").Value = 1748.1333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1748.1333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1546.1333
$ws.Range("N61").Value = $null
$ws.Range("H63").Value = 50750
$ws.Range("J63").Value = 50750
$ws.Range("L63").Value = 50750
$ws.Range("N63").Value = -52248
$ws.Range("H66").Value = 50750
$ws.Range("J66").Value = 50750
$ws.Range("L66").Value = 152250
$ws.Range("N66").Value = -159738
$ws.Range("H113").Value = 1748.1333
$ws.Range("I113").Value = 1748.1333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1748.1333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 421.8667
$ws.Range("N113").Value = $null
$ws.Range("H122").Value = 4984.9287
$ws.Range("I122").Value = 4286.5
$ws.Range("K122").Value = 12859.5
$ws.Range("M122").Value = -10409.5
$ws.Range("H136").Value = 29413596
$ws.Range("I136").Value = 1601.7858
$ws.Range("J136").Value = 166669570
$ws.Range("K136").Value = 4805.357400000001
$ws.Range("L136").Value = 500008710
$ws.Range("M136").Value = -2255.357400000001
$ws.Range("N136").Value = -500013810
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4381.923
$ws.Range("I81").Value = 2618
$ws.Range("J81").Value = 5484.375
$ws.Range("K81").Value = 5236
$ws.Range("L81").Value = 10968.75
$ws.Range("M81").Value = -4175
$ws.Range("N81").Value = -13090.75
$ws.Range("H84").Value = 4381.923
$ws.Range("I84").Value = 2618
$ws.Range("J84").Value = 5484.375
$ws.Range("K84").Value = 26180
$ws.Range("L84").Value = 54843.75
$ws.Range("M84").Value = -20876
$ws.Range("N84").Value = -65451.75
$ws.Range("H113").Value = 394
$ws.Range("I113").Value = 394
$ws.Range("K113").Value = 1182
$ws.Range("M113").Value = 988
$ws.Range("H122").Value = 5461.7334
$ws.Range("I122").Value = 5479.077
$ws.Range("J122").Value = 5349
$ws.Range("K122").Value = 16437.231
$ws.Range("L122").Value = 16047
$ws.Range("M122").Value = -13987.231
$ws.Range("N122").Value = -20947
$ws.Range("H136").Value = 7292.04
$ws.Range("J136").Value = 9111.764999999999
$ws.Range("L136").Value = 27335.295
$ws.Range("N136").Value = -32435.295

Write-Host "Applied Phantom_Profits updates"